$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the dates in column D between rows 2-3 and rows 4-5.
$ws.Range("D2").Value = 44559
$ws.Range("D3").Value = 44559
$ws.Range("D4").Value = 44574
$ws.Range("D5").Value = 44574
